$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 244.96552
$ws.Range("I55").Value = 264.94446
$ws.Range("K55").Value = 264.94446
$ws.Range("M55").Value = -50.94445999999999

$ws.Range("H98").Value = 7028.9
$ws.Range("I98").Value = 7311.125
$ws.Range("J98").Value = 5900
$ws.Range("K98").Value = 7311.125
$ws.Range("L98").Value = 5900
$ws.Range("M98").Value = -5813.125
$ws.Range("N98").Value = -8896

$ws.Range("H122").Value = 7028.9
$ws.Range("I122").Value = 7311.125
$ws.Range("J122").Value = 5900
$ws.Range("K122").Value = 21933.375
$ws.Range("L122").Value = 17700
$ws.Range("M122").Value = -19483.375
$ws.Range("N122").Value = -22600

$ws.Range("H132").Value = 1681.3043
$ws.Range("I132").Value = 1591.091
$ws.Range("K132").Value = 4773.272999999999
$ws.Range("M132").Value = -2243.272999999999

$ws.Range("H137").Value = 5001.6
$ws.Range("I137").Value = 5150.933
$ws.Range("K137").Value = 15452.799
$ws.Range("M137").Value = -12902.799

$ws.Range("H138").Value = 2003889.8
$ws.Range("J138").Value = 3129974
$ws.Range("L138").Value = 9389922
$ws.Range("N138").Value = -9400202

$ws.Range("H141").Value = 4670.737
$ws.Range("I141").Value = 4616.3335
$ws.Range("K141").Value = 13849.0005
$ws.Range("M141").Value = -8669.000499999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2197247
$ws.Range("I32").Value = 2845677.8
$ws.Range("K32").Value = 2845677.8
$ws.Range("M32").Value = -2845390.8

$ws.Range("H61").Value = 5200.8247
$ws.Range("I61").Value = 2584.5112
$ws.Range("J61").Value = 15012
$ws.Range("K61").Value = 2584.5112
$ws.Range("L61").Value = 15012
$ws.Range("M61").Value = -2372.5112
$ws.Range("N61").Value = -15436

$ws.Range("H110").Value = 22223588
$ws.Range("I110").Value = 1108.3334
$ws.Range("J110").Value = 55557308
$ws.Range("K110").Value = 1108.3334
$ws.Range("L110").Value = 55557308
$ws.Range("M110").Value = 936.6666
$ws.Range("N110").Value = -55561398

$ws.Range("H122").Value = 10524.207
$ws.Range("I122").Value = 15893.0625
$ws.Range("J122").Value = 3916.3845
$ws.Range("K122").Value = 47679.1875
$ws.Range("L122").Value = 11749.1535
$ws.Range("M122").Value = -45229.1875
$ws.Range("N122").Value = -16649.1535

$ws.Range("H132").Value = 3301487
$ws.Range("I132").Value = 10544451
$ws.Range("J132").Value = 9230.909
$ws.Range("K132").Value = 31633353
$ws.Range("L132").Value = 27692.727
$ws.Range("M132").Value = -31630823
$ws.Range("N132").Value = -32752.727

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H135").Value = 20000000
$ws.Range("J135").Value = 20000000
$ws.Range("L135").Value = 20000000
$ws.Range("N135").Value = -20010140

$ws.Range("H136").Value = 5200.8247
$ws.Range("I136").Value = 2584.5112
$ws.Range("J136").Value = 15012
$ws.Range("K136").Value = 7753.5336
$ws.Range("L136").Value = 45036
$ws.Range("M136").Value = -5203.5336
$ws.Range("N136").Value = -50136

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5559.074
$ws.Range("I134").Value = 1689.4736
$ws.Range("J134").Value = 14749.375
$ws.Range("K134").Value = 5068.4208
$ws.Range("L134").Value = 44248.125
$ws.Range("M134").Value = -2533.4208
$ws.Range("N134").Value = -49318.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6952.769
$ws.Range("I16").Value = 6929.625
$ws.Range("J16").Value = 6989.8
$ws.Range("K16").Value = 6929.625
$ws.Range("L16").Value = 6989.8
$ws.Range("M16").Value = -6642.625
$ws.Range("N16").Value = -7563.8

$ws.Range("H31").Value = 5634.089
$ws.Range("I31").Value = 2554.9167
$ws.Range("J31").Value = 9153.143
$ws.Range("K31").Value = 2554.9167
$ws.Range("L31").Value = 9153.143
$ws.Range("M31").Value = -2259.9167
$ws.Range("N31").Value = -9743.143

$ws.Range("H34").Value = 5634.089
$ws.Range("I34").Value = 2554.9167
$ws.Range("J34").Value = 9153.143
$ws.Range("K34").Value = 2554.9167
$ws.Range("L34").Value = 9153.143
$ws.Range("M34").Value = -2352.9167
$ws.Range("N34").Value = -9557.143

$ws.Range("H99").Value = 7163.8887
$ws.Range("I99").Value = 7187.6665
$ws.Range("K99").Value = 7187.6665
$ws.Range("M99").Value = -5689.6665

$ws.Range("H107").Value = 2437.8262
$ws.Range("I107").Value = 1985.6875
$ws.Range("J107").Value = 3471.2856
$ws.Range("K107").Value = 1985.6875
$ws.Range("L107").Value = 3471.2856
$ws.Range("M107").Value = -65.6875
$ws.Range("N107").Value = -7311.2856

$ws.Range("H113").Value = 6952.769
$ws.Range("I113").Value = 6929.625
$ws.Range("J113").Value = 6989.8
$ws.Range("K113").Value = 6929.625
$ws.Range("L113").Value = 6989.8
$ws.Range("M113").Value = -4759.625
$ws.Range("N113").Value = -11329.8

$ws.Range("H126").Value = 7163.8887
$ws.Range("I126").Value = 7187.6665
$ws.Range("K126").Value = 21562.9995
$ws.Range("M126").Value = -19092.9995

$ws.Range("H132").Value = 8057.75
$ws.Range("I132").Value = 2956.8572
$ws.Range("K132").Value = 8870.571599999999
$ws.Range("M132").Value = -6340.571599999999

$ws.Range("H134").Value = 7607.3076
$ws.Range("I134").Value = 1737.5
$ws.Range("J134").Value = 16999
$ws.Range("K134").Value = 5212.5
$ws.Range("L134").Value = 50997
$ws.Range("M134").Value = -2677.5
$ws.Range("N134").Value = -56067

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I59").Value = 2000
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 6000
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -5460
$ws.Range("N59").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 105
$ws.Range("I9").Value = 105
$ws.Range("K9").Value = 105
$ws.Range("M9").Value = 65

$ws.Range("H70").Value = 10798.667
$ws.Range("I70").Value = 8396.5
$ws.Range("J70").Value = 11999.75
$ws.Range("K70").Value = 8396.5
$ws.Range("L70").Value = 11999.75
$ws.Range("M70").Value = -8126.5
$ws.Range("N70").Value = -12539.75

$ws.Range("H73").Value = 10798.667
$ws.Range("I73").Value = 8396.5
$ws.Range("J73").Value = 11999.75
$ws.Range("K73").Value = 8396.5
$ws.Range("L73").Value = 11999.75
$ws.Range("M73").Value = -7460.5
$ws.Range("N73").Value = -13871.75

$ws.Range("H102").Value = 1383.6818
$ws.Range("I102").Value = 1383.6818
$ws.Range("K102").Value = 1383.6818
$ws.Range("M102").Value = 238.3181999999999

$ws.Range("H122").Value = 125131230
$ws.Range("I122").Value = 333669340
$ws.Range("J122").Value = 8379
$ws.Range("K122").Value = 1001008020
$ws.Range("L122").Value = 25137
$ws.Range("M122").Value = -1001005570
$ws.Range("N122").Value = -30037

$ws.Range("H132").Value = 4998.4443
$ws.Range("I132").Value = 1823.4286
$ws.Range("J132").Value = 16111
$ws.Range("K132").Value = 5470.2858
$ws.Range("L132").Value = 48333
$ws.Range("M132").Value = -2940.2858
$ws.Range("N132").Value = -53393

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7108.0713
$ws.Range("I132").Value = 3158.2354
$ws.Range("K132").Value = 9474.706200000001
$ws.Range("M132").Value = -6944.706200000001

$ws.Range("H136").Value = 15909.643
$ws.Range("I136").Value = 9073.6
$ws.Range("J136").Value = 32999.75
$ws.Range("K136").Value = 27220.8
$ws.Range("L136").Value = 98999.25
$ws.Range("M136").Value = -24670.8
$ws.Range("N136").Value = -104099.25

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws.Range("H141").Value = 89807
$ws.Range("J141").Value = 89807
$ws.Range("L141").Value = 89807
$ws.Range("N141").Value = -100167

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 46479.54
$ws.Range("I132").Value = 52650.332
$ws.Range("J132").Value = 41190.285
$ws.Range("K132").Value = 157950.996
$ws.Range("L132").Value = 123570.855
$ws.Range("M132").Value = -155420.996
$ws.Range("N132").Value = -128630.855
